# Updates the coin/price/volume table to match the refreshed crypto snapshot.
# All Price/Volume(1h) values are stored as text in the sheet (t="inlineStr"/shared-string),
# so every numeric-looking value is written with a leading apostrophe to force Excel to
# keep it as text instead of silently coercing it to a Number/Percentage cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''245.70'
$ws.Range("E2").Value = '''-0.39%'

# Row 3
$ws.Range("D3").Value = '''30.09'
$ws.Range("E3").Value = '''-0.93%'

# Row 4
$ws.Range("D4").Value = '''5.157'
$ws.Range("E4").Value = '''-0.47%'

# Row 5
$ws.Range("D5").Value = '''0.05764'
$ws.Range("E5").Value = '''0.58%'

# Row 6
$ws.Range("D6").Value = '''6.664'
$ws.Range("E6").Value = '''0.95%'

# Row 7
$ws.Range("D7").Value = '''3.283'
$ws.Range("E7").Value = '''6.86%'

# Row 8
$ws.Range("D8").Value = '''0.8494'
$ws.Range("E8").Value = '''-0.72%'

# Row 9
$ws.Range("D9").Value = '''0.8590'
$ws.Range("E9").Value = '''-2.61%'

# Row 10
$ws.Range("D10").Value = '''0.1387'
$ws.Range("E10").Value = '''1.61%'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03417'
$ws.Range("E11").Value = '''5.17%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07078'
$ws.Range("E12").Value = '''-0.46%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03235'
$ws.Range("E13").Value = '''12.68%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09361'
$ws.Range("E14").Value = '''-0.35%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001526'
$ws.Range("E15").Value = '''0.17%'

# Row 16
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '''0.0005937'
$ws.Range("E16").Value = '''-1.58%'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.005949'
$ws.Range("E17").Value = '''-0.26%'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.548'
$ws.Range("E18").Value = '''1.64%'

# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.201'
$ws.Range("E19").Value = '''-3.10%'

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3145'
$ws.Range("E20").Value = '''-2.42%'

# Row 21
$ws.Range("D21").Value = '''0.1315'
$ws.Range("E21").Value = '''1.12%'

# Row 22
$ws.Range("D22").Value = '''3.479'

# Row 23
$ws.Range("E23").Value = '''2.11%'

# Row 24
$ws.Range("E24").Value = '''-1.16%'

# Row 25
$ws.Range("D25").Value = '''0.001222'
$ws.Range("E25").Value = '''0.79%'

# Row 26
$ws.Range("D26").Value = '''0.004161'
$ws.Range("E26").Value = '''-7.31%'

# Row 27
$ws.Range("E27").Value = '''-0.81%'

# Row 28
$ws.Range("D28").Value = '''0.0001448'
$ws.Range("E28").Value = '''4.66%'

# Row 40
$ws.Range("E40").Value = '''-0.77%'

# Row 41
$ws.Range("D41").Value = '''0.1071'
$ws.Range("E41").Value = '''0.17%'

# Row 42
$ws.Range("D42").Value = '''0.002460'
$ws.Range("E42").Value = '''11.85%'

# Row 43
$ws.Range("E43").Value = '''-48.32%'

# Row 44
$ws.Range("D44").Value = '''0.01022'
$ws.Range("E44").Value = '''2.36%'

# Row 45
$ws.Range("E45").Value = '''7.66%'

# Row 46
$ws.Range("E46").Value = '''0.00%'

# Row 47
$ws.Range("D47").Value = '''0.07097'
$ws.Range("E47").Value = '''-11.33%'

# Row 48
$ws.Range("D48").Value = '''0.002473'
$ws.Range("E48").Value = '''-10.53%'

# Row 49
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("E49").Value = '''0.00%'

# Row 50
$ws.Range("D50").Value = '''0.0001999'
$ws.Range("E50").Value = '''0.00%'
